# ---------------------------------------------------------------------------
# Re-organise the workbook:
#   * insert two new leading worksheets (TXL_SHEET_TRACKER, TEST_SHEET_TXL_ROW)
#   * populate them
#   * populate the pre-existing TEST_STANDARD_ROW sheet with index/color rows
#   * refresh cursor/selection on every sheet
#   * replace the single "YO" defined name with four new named ranges
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. insert the two new worksheets in front of TEST_SHEET -------------
# sheetId numbering is assigned at creation time (not at final position), so
# create TXL_SHEET_TRACKER first (-> sheetId 10) and TEST_SHEET_TXL_ROW
# second (-> sheetId 11), then slot TEST_SHEET_TXL_ROW in front of the
# tracker sheet so the final tab order is:
#   TEST_SHEET_TXL_ROW, TXL_SHEET_TRACKER, TEST_SHEET, TEST_STANDARD_ROW, TEST_SCALAR_INPUT
$testSheet = $wb.Worksheets.Item("TEST_SHEET")

$tracker = $wb.Worksheets.Add($testSheet)
$tracker.Name = "TXL_SHEET_TRACKER"

$txlRow = $wb.Worksheets.Add($tracker)
$txlRow.Name = "TEST_SHEET_TXL_ROW"

# Worksheet handles obtained before an Add() track a *slot index*, not a
# stable sheet identity, so every handle grabbed earlier now resolves to the
# wrong sheet (it "slid" along with the insertions). Re-resolve every sheet
# we need to touch, by name, now that the tab order is final.
$txlRow      = $wb.Worksheets.Item("TEST_SHEET_TXL_ROW")
$tracker     = $wb.Worksheets.Item("TXL_SHEET_TRACKER")
$testSheet   = $wb.Worksheets.Item("TEST_SHEET")
$standardRow = $wb.Worksheets.Item("TEST_STANDARD_ROW")
$scalarInput = $wb.Worksheets.Item("TEST_SCALAR_INPUT")

# --- 2. populate TXL_SHEET_TRACKER ----------------------------------------
$tracker.Range("A1").Value = "sheet_name"
$tracker.Range("B1").Value = "descr"
$tracker.Range("C1").Value = "sheet_type"
$tracker.Range("A2").Value = "TEST_SCALAR_INPUT"
$tracker.Range("B2").Value = "A Test worksheet with scalar inputs"
$tracker.Range("C2").Value = 0

$tracker.Range("A1:C1").Interior.Color = 0xBFBFBF

$tracker.Columns.Item(1).ColumnWidth = 20.375
$tracker.Columns.Item(2).ColumnWidth = 30.25
$tracker.Columns.Item(3).ColumnWidth = 9.75
$tracker.Columns.Item(4).ColumnWidth = 1.875
$tracker.Columns.Item(5).ColumnWidth = 1.875

# --- 3. populate TEST_STANDARD_ROW ----------------------------------------
$standardRow.Range("A1").Value = "index"
$standardRow.Range("B1").Value = "color"
$standardRow.Range("A2").Value = 1
$standardRow.Range("B2").Value = "blue"
$standardRow.Range("A3").Value = 2
$standardRow.Range("B3").Value = "red"
$standardRow.Range("A4").Value = 3
$standardRow.Range("B4").Value = "green"
$standardRow.Range("A5").Value = 4
$standardRow.Range("B5").Value = "yellow"

# --- 4. refresh the selection / active cell on every sheet ----------------
$txlRow.Range("E38").Select()
$tracker.Range("B23").Select()
$testSheet.Range("H38").Select()
$scalarInput.Range("B3").Select()

# TEST_STANDARD_ROW ends up the active tab
$standardRow.Activate()
$standardRow.Range("B3").Select()

# --- 5. defined names: drop "YO", add the four new names ------------------
$wb.Names.Item("YO").Delete()

$wb.Names.Add('TEST_NAME_FOR_NBK', '=TEST_SHEET_TXL_ROW!$I$13')
$wb.Names.Add('TEST_SCALAR_INPUT__dollar', '=TEST_SCALAR_INPUT!$B$2')
$wb.Names.Add('TEST_SCALAR_INPUT__kwh', '=TEST_SCALAR_INPUT!$B$3')
$wb.Names.Add('TESTOTHER_NAME', '=TEST_SHEET_TXL_ROW!$O$12')
